$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 9808.083000000001
$ws.Range("I64").Value = 7969.8
$ws.Range("K64").Value = 7969.8
$ws.Range("M64").Value = -7721.8

$ws.Range("H67").Value = 9808.083000000001
$ws.Range("I67").Value = 7969.8
$ws.Range("K67").Value = 7969.8
$ws.Range("M67").Value = -7111.8

$ws.Range("H116").Value = 6496.613
$ws.Range("I116").Value = 6241.3335
$ws.Range("K116").Value = 6241.3335
$ws.Range("M116").Value = -2799.3335

$ws.Range("H138").Value = 7003.027
$ws.Range("J138").Value = 8296.218000000001
$ws.Range("L138").Value = 24888.654
$ws.Range("N138").Value = -35168.654

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H61").Value = 2432.25
$ws.Range("I61").Value = 2286.1428
$ws.Range("K61").Value = 2286.1428
$ws.Range("M61").Value = -2074.1428

$ws.Range("H88").Value = 5134.0586
$ws.Range("I88").Value = 5068.143
$ws.Range("J88").Value = 5180.2
$ws.Range("K88").Value = 5068.143
$ws.Range("L88").Value = 5180.2
$ws.Range("M88").Value = -4662.143
$ws.Range("N88").Value = -5992.2

$ws.Range("H91").Value = 5134.0586
$ws.Range("I91").Value = 5068.143
$ws.Range("J91").Value = 5180.2
$ws.Range("K91").Value = 5068.143
$ws.Range("L91").Value = 5180.2
$ws.Range("M91").Value = -3664.143
$ws.Range("N91").Value = -7988.2

$ws.Range("H110").Value = 2261.25
$ws.Range("I110").Value = 1381.6364
$ws.Range("K110").Value = 1381.6364
$ws.Range("M110").Value = 663.3635999999999

$ws.Range("H122").Value = 1436
$ws.Range("I122").Value = 1394.2727
$ws.Range("K122").Value = 4182.8181
$ws.Range("M122").Value = -1732.8181

$ws.Range("H136").Value = 2432.25
$ws.Range("I136").Value = 2286.1428
$ws.Range("K136").Value = 6858.428400000001
$ws.Range("M136").Value = -4308.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3053
$ws.Range("I86").Value = 2380
$ws.Range("K86").Value = 2380
$ws.Range("M86").Value = -1257

$ws.Range("H89").Value = 3053
$ws.Range("I89").Value = 2380
$ws.Range("K89").Value = 11900
$ws.Range("M89").Value = -6284

$ws.Range("H99").Value = 2686
$ws.Range("I99").Value = 2686
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2686
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1188
$ws.Range("N99").ClearContents()

$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -43134

$ws.Range("H134").Value = 2274.0908
$ws.Range("I134").Value = 1963.3334
$ws.Range("K134").Value = 5890.0002
$ws.Range("M134").Value = -3355.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 10028.75
$ws.Range("I41").Value = 3371.6667
$ws.Range("K41").Value = 3371.6667
$ws.Range("M41").Value = -2943.6667

$ws.Range("H105").Value = 3394.8333
$ws.Range("I105").Value = 991.1429000000001
$ws.Range("K105").Value = 991.1429000000001
$ws.Range("M105").Value = 755.8570999999999

$ws.Range("H122").Value = 599.6
$ws.Range("J122").Value = 749.5
$ws.Range("L122").Value = 2248.5
$ws.Range("N122").Value = -7148.5

$ws.Range("H134").Value = 3504.0908
$ws.Range("I134").Value = 3054.7
$ws.Range("K134").Value = 9164.099999999999
$ws.Range("M134").Value = -6629.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 187848.75
$ws.Range("I12").Value = 333372.66
$ws.Range("J12").Value = 746.5714
$ws.Range("K12").Value = 1000117.98
$ws.Range("L12").Value = 2239.7142
$ws.Range("M12").Value = -999944.98
$ws.Range("N12").Value = -2585.7142

$ws.Range("H34").Value = 490.625
$ws.Range("J34").Value = 1025
$ws.Range("L34").Value = 3075
$ws.Range("N34").Value = -3243

$ws.Range("H39").Value = 2251.25
$ws.Range("I39").Value = 403
$ws.Range("J39").Value = 2515.2856
$ws.Range("K39").Value = 1209
$ws.Range("L39").Value = 7545.8568
$ws.Range("M39").Value = -915
$ws.Range("N39").Value = -8133.8568

$ws.Range("H55").Value = 2062
$ws.Range("J55").Value = 2582.8333
$ws.Range("L55").Value = 7748.499899999999
$ws.Range("N55").Value = -8102.499899999999

$ws.Range("H80").Value = 8482.666999999999
$ws.Range("I80").Value = 1966.3334
$ws.Range("J80").Value = 14999
$ws.Range("K80").Value = 5899.0002
$ws.Range("L80").Value = 44997
$ws.Range("M80").Value = -4963.0002
$ws.Range("N80").Value = -46869

$ws.Range("H83").Value = 8482.666999999999
$ws.Range("I83").Value = 1966.3334
$ws.Range("J83").Value = 14999
$ws.Range("K83").Value = 17697.0006
$ws.Range("L83").Value = 134991
$ws.Range("M83").Value = -13017.0006
$ws.Range("N83").Value = -144351

$ws.Range("H92").Value = 775.25
$ws.Range("I92").Value = 434
$ws.Range("K92").Value = 1302
$ws.Range("M92").Value = -54

$ws.Range("H103").Value = 27122.9
$ws.Range("I103").Value = 83975
$ws.Range("J103").Value = 2757.7144
$ws.Range("K103").Value = 251925
$ws.Range("L103").Value = 8273.143199999999
$ws.Range("M103").Value = -251046
$ws.Range("N103").Value = -10031.1432

$ws.Range("H113").Value = 567.3333
$ws.Range("J113").Value = 800
$ws.Range("L113").Value = 2400
$ws.Range("N113").Value = -6740

$ws.Range("H132").Value = 2976.2307
$ws.Range("I132").Value = 1554.7142
$ws.Range("K132").Value = 13992.4278
$ws.Range("M132").Value = -11462.4278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2672.647
$ws.Range("I80").Value = 2419.3
$ws.Range("J80").Value = 3034.5715
$ws.Range("K80").Value = 2419.3
$ws.Range("L80").Value = 3034.5715
$ws.Range("M80").Value = -1421.3
$ws.Range("N80").Value = -5030.5715

$ws.Range("H83").Value = 2672.647
$ws.Range("I83").Value = 2419.3
$ws.Range("J83").Value = 3034.5715
$ws.Range("K83").Value = 12096.5
$ws.Range("L83").Value = 15172.8575
$ws.Range("M83").Value = -7104.5
$ws.Range("N83").Value = -25156.8575

$ws.Range("H102").Value = 493.75
$ws.Range("I102").Value = 500
$ws.Range("J102").Value = 475
$ws.Range("K102").Value = 500
$ws.Range("L102").Value = 475
$ws.Range("M102").Value = 1122
$ws.Range("N102").Value = -3719

$ws.Range("H122").Value = 85955.664
$ws.Range("I122").Value = 102349.3
$ws.Range("K122").Value = 307047.9
$ws.Range("M122").Value = -304597.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -888

$ws.Range("H68").Value = 2175
$ws.Range("I68").Value = 2175
$ws.Range("K68").Value = 2175
$ws.Range("M68").Value = -1426

$ws.Range("H71").Value = 2175
$ws.Range("I71").Value = 2175
$ws.Range("K71").Value = 10875
$ws.Range("M71").Value = -7131

$ws.Range("H74").Value = 39999.668
$ws.Range("J74").Value = 39999.668
$ws.Range("L74").Value = 39999.668
$ws.Range("N74").Value = -41995.668

$ws.Range("H77").Value = 39999.668
$ws.Range("J77").Value = 39999.668
$ws.Range("L77").Value = 119999.004
$ws.Range("N77").Value = -129983.004

$ws.Range("H82").Value = 2625.0356
$ws.Range("I82").Value = 1934.8
$ws.Range("J82").Value = 3008.5
$ws.Range("K82").Value = 1934.8
$ws.Range("L82").Value = 3008.5
$ws.Range("M82").Value = -1573.8
$ws.Range("N82").Value = -3730.5

$ws.Range("H85").Value = 2625.0356
$ws.Range("I85").Value = 1934.8
$ws.Range("J85").Value = 3008.5
$ws.Range("K85").Value = 1934.8
$ws.Range("L85").Value = 3008.5
$ws.Range("M85").Value = -686.8
$ws.Range("N85").Value = -5504.5

$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2152.16
$ws.Range("I132").Value = 2106.7273
$ws.Range("K132").Value = 6320.1819
$ws.Range("M132").Value = -3790.1819
